$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2023-05-05 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-05-06 Saturday", 2) | Out-Null

# Update each cell of the multiplication-answers table by explicit
# (row, column) address so duplicate cell text (e.g. two cells both
# originally "57×37=2109") are updated independently and correctly.
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "10×74=740"
$t.Cell(1,2).Range.Text = "55×65=3575"
$t.Cell(1,3).Range.Text = "86×56=4816"
$t.Cell(1,4).Range.Text = "65×76=4940"
$t.Cell(1,5).Range.Text = "58×59=3422"

$t.Cell(2,1).Range.Text = "29×89=2581"
$t.Cell(2,2).Range.Text = "99×76=7524"
$t.Cell(2,3).Range.Text = "35×41=1435"
$t.Cell(2,4).Range.Text = "71×69=4899"
$t.Cell(2,5).Range.Text = "96×12=1152"

$t.Cell(3,1).Range.Text = "54×55=2970"
$t.Cell(3,2).Range.Text = "26×12=312"
$t.Cell(3,3).Range.Text = "46×95=4370"
$t.Cell(3,4).Range.Text = "40×99=3960"
$t.Cell(3,5).Range.Text = "92×89=8188"

$t.Cell(4,1).Range.Text = "42×20=840"
$t.Cell(4,2).Range.Text = "27×37=999"
$t.Cell(4,3).Range.Text = "95×37=3515"
$t.Cell(4,4).Range.Text = "38×48=1824"
$t.Cell(4,5).Range.Text = "43×60=2580"

$t.Cell(5,1).Range.Text = "15×13=195"
$t.Cell(5,2).Range.Text = "95×68=6460"
$t.Cell(5,3).Range.Text = "38×46=1748"
$t.Cell(5,4).Range.Text = "68×89=6052"
$t.Cell(5,5).Range.Text = "37×45=1665"

$t.Cell(6,1).Range.Text = "12×62=744"
$t.Cell(6,2).Range.Text = "80×16=1280"
$t.Cell(6,3).Range.Text = "76×50=3800"
$t.Cell(6,4).Range.Text = "10×94=940"
$t.Cell(6,5).Range.Text = "65×57=3705"

$t.Cell(7,1).Range.Text = "67×72=4824"
$t.Cell(7,2).Range.Text = "35×33=1155"
$t.Cell(7,3).Range.Text = "31×29=899"
$t.Cell(7,4).Range.Text = "99×91=9009"
$t.Cell(7,5).Range.Text = "76×24=1824"

$t.Cell(8,1).Range.Text = "47×44=2068"
$t.Cell(8,2).Range.Text = "18×17=306"
$t.Cell(8,3).Range.Text = "16×11=176"
$t.Cell(8,4).Range.Text = "37×79=2923"
$t.Cell(8,5).Range.Text = "79×64=5056"

$t.Cell(9,1).Range.Text = "90×31=2790"
$t.Cell(9,2).Range.Text = "68×49=3332"
$t.Cell(9,3).Range.Text = "27×88=2376"
$t.Cell(9,4).Range.Text = "45×58=2610"
$t.Cell(9,5).Range.Text = "63×54=3402"

$t.Cell(10,1).Range.Text = "59×89=5251"
$t.Cell(10,2).Range.Text = "30×68=2040"
$t.Cell(10,3).Range.Text = "38×85=3230"
$t.Cell(10,4).Range.Text = "82×45=3690"
$t.Cell(10,5).Range.Text = "52×80=4160"

$t.Cell(11,1).Range.Text = "57×26=1482"
$t.Cell(11,2).Range.Text = "82×95=7790"
$t.Cell(11,3).Range.Text = "58×54=3132"
$t.Cell(11,4).Range.Text = "10×38=380"
$t.Cell(11,5).Range.Text = "66×90=5940"

$t.Cell(12,1).Range.Text = "92×85=7820"
$t.Cell(12,2).Range.Text = "65×87=5655"
$t.Cell(12,3).Range.Text = "19×46=874"
$t.Cell(12,4).Range.Text = "56×76=4256"
$t.Cell(12,5).Range.Text = "45×95=4275"

$t.Cell(13,1).Range.Text = "51×85=4335"
$t.Cell(13,2).Range.Text = "30×12=360"
$t.Cell(13,3).Range.Text = "15×44=660"
$t.Cell(13,4).Range.Text = "25×57=1425"
$t.Cell(13,5).Range.Text = "18×30=540"

$t.Cell(14,1).Range.Text = "29×19=551"
$t.Cell(14,2).Range.Text = "33×72=2376"
$t.Cell(14,3).Range.Text = "67×38=2546"
$t.Cell(14,4).Range.Text = "34×38=1292"
$t.Cell(14,5).Range.Text = "80×26=2080"

$t.Cell(15,1).Range.Text = "31×47=1457"
$t.Cell(15,2).Range.Text = "14×77=1078"
$t.Cell(15,3).Range.Text = "42×68=2856"
$t.Cell(15,4).Range.Text = "28×35=980"
$t.Cell(15,5).Range.Text = "47×90=4230"

$t.Cell(16,1).Range.Text = "29×91=2639"
$t.Cell(16,2).Range.Text = "42×81=3402"
$t.Cell(16,3).Range.Text = "48×84=4032"
$t.Cell(16,4).Range.Text = "45×13=585"
$t.Cell(16,5).Range.Text = "66×85=5610"

$t.Cell(17,1).Range.Text = "59×75=4425"
$t.Cell(17,2).Range.Text = "37×76=2812"
$t.Cell(17,3).Range.Text = "97×43=4171"
$t.Cell(17,4).Range.Text = "99×81=8019"
$t.Cell(17,5).Range.Text = "76×66=5016"

$t.Cell(18,1).Range.Text = "32×97=3104"
$t.Cell(18,2).Range.Text = "27×74=1998"
$t.Cell(18,3).Range.Text = "88×84=7392"
$t.Cell(18,4).Range.Text = "68×35=2380"
$t.Cell(18,5).Range.Text = "46×34=1564"

$t.Cell(19,1).Range.Text = "41×36=1476"
$t.Cell(19,2).Range.Text = "10×31=310"
$t.Cell(19,3).Range.Text = "85×36=3060"
$t.Cell(19,4).Range.Text = "58×25=1450"
$t.Cell(19,5).Range.Text = "37×84=3108"

$t.Cell(20,1).Range.Text = "79×80=6320"
$t.Cell(20,2).Range.Text = "48×89=4272"
$t.Cell(20,3).Range.Text = "96×53=5088"
$t.Cell(20,4).Range.Text = "76×37=2812"
$t.Cell(20,5).Range.Text = "37×37=1369"

